$wb = $excel.ActiveWorkbook

# Sheet 1: 个人持仓 (Personal Holdings)
$ws1 = $wb.Worksheets.Item("个人持仓")
$ws1.Range("C2").Value = 52.41
$ws1.Range("C3").Value = 38.9
$ws1.Range("C4").Value = 49.9
$ws1.Range("C5").Value = 53.83
$ws1.Range("C6").Value = 31.65
$ws1.Range("C7").Value = 41.44
$ws1.Range("C8").Value = 25.17
$ws1.Range("C9").Value = 24.25
$ws1.Range("C10").Value = 25.58
$ws1.Range("C11").Value = 144.6
$ws1.Range("C12").Value = 178.6
$ws1.Range("C13").Value = 237.88
$ws1.Range("C14").Value = 0.847
$ws1.Range("C15").Value = 11.33
$ws1.Range("C16").Value = 28.51
$ws1.Range("C17").Value = 25.33
$ws1.Range("C18").Value = 18.89
$ws1.Range("C19").Value = 37.1
$ws1.Range("C20").Value = 38.22
$ws1.Range("C21").Value = 28.75
$ws1.Range("C22").Value = 124.13
$ws1.Range("C23").Value = 4.073
$ws1.Range("C24").Value = 4.406

# Sheet 2: 家庭持仓 (Family Holdings)
$ws2 = $wb.Worksheets.Item("家庭持仓")
$ws2.Range("C2").Value = 52.41
$ws2.Range("C3").Value = 38.9
$ws2.Range("C4").Value = 49.9
$ws2.Range("C5").Value = 53.83
$ws2.Range("C6").Value = 31.65
$ws2.Range("C7").Value = 41.44
$ws2.Range("C8").Value = 25.17
$ws2.Range("C9").Value = 24.25
$ws2.Range("C10").Value = 25.58
$ws2.Range("C11").Value = 144.6
$ws2.Range("C12").Value = 237.88
$ws2.Range("C13").Value = 178.6
$ws2.Range("C14").Value = 234.57
$ws2.Range("C15").Value = 0.847
$ws2.Range("C16").Value = 46.31
$ws2.Range("C17").Value = 11.33
$ws2.Range("C18").Value = 28.51
$ws2.Range("C19").Value = 25.33
$ws2.Range("C20").Value = 18.89
$ws2.Range("C21").Value = 37.1
$ws2.Range("C22").Value = 38.22
$ws2.Range("C23").Value = 28.75
$ws2.Range("C24").Value = 124.13
$ws2.Range("C25").Value = 4.073
$ws2.Range("C26").Value = 1.17
$ws2.Range("C27").Value = 4.406
